$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 3 (Auxerre)
$ws.Range("C3").Value = 25.5

# Row 4 (Brest)
$ws.Range("D4").Value = 42.7
$ws.Range("E4").Value = 23
$ws.Range("F4").Value = 253
$ws.Range("G4").Value = 2070
$ws.Range("H4").Value = 23
$ws.Range("I4").Value = 30
$ws.Range("J4").Value = 19
$ws.Range("K4").Value = 49
$ws.Range("L4").Value = 24
$ws.Range("O4").Value = 40
$ws.Range("Q4").Value = 1.3
$ws.Range("R4").Value = 0.83
$ws.Range("S4").Value = 2.13
$ws.Range("T4").Value = 1.04
$ws.Range("U4").Value = 1.87

# Row 10 (Marseille)
$ws.Range("C10").Value = 28.1
$ws.Range("D10").Value = 58.2
$ws.Range("E10").Value = 23
$ws.Range("F10").Value = 253
$ws.Range("G10").Value = 2070
$ws.Range("H10").Value = 23
$ws.Range("N10").Value = 7
$ws.Range("O10").Value = 51
$ws.Range("Q10").Value = 2.04
$ws.Range("R10").Value = 1.52
$ws.Range("S10").Value = 3.57
$ws.Range("T10").Value = 1.83
$ws.Range("U10").Value = 3.35
